$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E6").Value = "Trabaja"
$ws.Range("E9").Value = "Nada"
$ws.Range("E11").Value = "Trabaja"
$ws.Range("E18").Value = "Nada"
$ws.Range("E19").Value = "Nada"
$ws.Range("E24").Value = "Trabaja"
$ws.Range("E25").Value = "Trabaja"
$ws.Range("E31").Value = "Nada"
$ws.Range("E34").Value = "Trabaja"
$ws.Range("E61").Value = "Trabaja"
$ws.Range("E64").Value = "Nada"
$ws.Range("E65").Value = "Nada"
$ws.Range("E70").Value = "Trabaja"
$ws.Range("E71").Value = "Trabaja"
$ws.Range("E77").Value = "Nada"
$ws.Range("E80").Value = "Trabaja"
$ws.Range("E94").Value = "Trabaja"
$ws.Range("E95").Value = "Trabaja"
$ws.Range("E121").Value = "Trabaja"
$ws.Range("E123").Value = "Nada"
$ws.Range("E126").Value = "Trabaja"
$ws.Range("E132").Value = "Nada"
$ws.Range("E133").Value = "Nada"
$ws.Range("E140").Value = "Trabaja"
$ws.Range("E141").Value = "Trabaja"
$ws.Range("E168").Value = "Trabaja"
$ws.Range("E170").Value = "Nada"
$ws.Range("E173").Value = "Trabaja"
$ws.Range("E178").Value = "Nada"
$ws.Range("E179").Value = "Nada"
$ws.Range("E186").Value = "Trabaja"
$ws.Range("E187").Value = "Trabaja"
$ws.Range("E193").Value = "Trabaja"
$ws.Range("E194").Value = "Nada"
$ws.Range("E202").Value = "Trabaja"
$ws.Range("E208").Value = "Nada"
$ws.Range("E213").Value = "Trabaja"
$ws.Range("E217").Value = "Nada"
$ws.Range("E218").Value = "Trabaja"
$ws.Range("E224").Value = "Nada"
$ws.Range("E225").Value = "Nada"
$ws.Range("E233").Value = "Trabaja"
$ws.Range("E239").Value = "Trabaja"
$ws.Range("E241").Value = "Nada"
$ws.Range("E248").Value = "Trabaja"
$ws.Range("E249").Value = "Trabaja"
$ws.Range("E254").Value = "Nada"
$ws.Range("E255").Value = "Nada"
$ws.Range("E262").Value = "Trabaja"
$ws.Range("E264").Value = "Nada"
$ws.Range("E267").Value = "Trabaja"
$ws.Range("E271").Value = "Nada"
$ws.Range("E285").Value = "Trabaja"
$ws.Range("E288").Value = "Nada"
$ws.Range("E294").Value = "Trabaja"
$ws.Range("E295").Value = "Trabaja"
$ws.Range("E296").Value = "Trabaja"
$ws.Range("E300").Value = "Nada"
$ws.Range("E301").Value = "Nada"
$ws.Range("E302").Value = "Nada"
$ws.Range("E304").Value = "Trabaja"
$ws.Range("E309").Value = "Trabaja"
$ws.Range("E310").Value = "Nada"
